$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 254, shifting existing rows 254:269 down to 255:270
$ws.Rows.Item(254).Insert()

# Populate the newly inserted row 254 with the new record
$ws.Range("A254").Value = 3
$ws.Range("B254").Value = "Femacal de La Calera"
$ws.Range("C254").Value = "Coquimbo"
$ws.Range("D254").Value = 44585
$ws.Range("E254").Value = 5
$ws.Range("F254").Value = 100112012
$ws.Range("G254").Value = "Espinaca"
$ws.Range("H254").Value = "Sin especificar"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 180
$ws.Range("K254").Value = 3500
$ws.Range("L254").Value = 4000
$ws.Range("M254").Value = 3750
$ws.Range("N254").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O254").Value = "Provincia de Quillota"
$ws.Range("P254").Value = 1250
$ws.Range("Q254").Value = 3
$ws.Range("R254").Value = "Hortaliza"
